# Updated cryptos list (price + 1h volume refresh) for rows 2-51 of Sheet1.
#
# Notes:
#  - Column D ("Price") values are stored as literal text in the workbook
#    (e.g. "1.002"), not numbers. Many of the new prices look like plain
#    decimals (e.g. "1.001"), which Excel would otherwise auto-convert to
#    a number on assignment. We prefix those with a leading apostrophe so
#    Excel keeps them as text, matching the original formatting. Values
#    that already can't parse as a number (multiple dots, subscript
#    digits, etc.) are assigned as-is.
#  - Rows 37/38 also swap which coin (MXToken / ImmutableX) occupies each
#    row, so Coin/Link/Price/Volume are all rewritten for those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.898.75"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "1.632.45"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'216.05"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "'0.5120"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.06339"
$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("D10").Value = "'19.47"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").Value = "'0.07769"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").Value = "'4.260"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").Value = "1.631.07"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").Value = "1.856.13"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").Value = "'0.5508"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").Value = "'63.76"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").Value = "0.0₅7639"
$ws.Range("E17").Value = "  -1.35%  "

$ws.Range("D18").Value = "25.916.17"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").Value = "'195.04"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").Value = "'4.417"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "'9.859"
$ws.Range("E22").Value = "  -0.72%  "

$ws.Range("D23").Value = "'6.026"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("D25").Value = "'1.891"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").Value = "'142.18"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").Value = "'0.1254"
$ws.Range("E27").Value = "  +4.82%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").Value = "'6.759"
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("D30").Value = "'1.242"
$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("D31").Value = "'0.04898"
$ws.Range("E31").Value = "  +0.54%  "

$ws.Range("D32").Value = "'3.232"
$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("D33").Value = "'3.183"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").Value = "'1.543"

$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").Value = "'0.8979"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.5499"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.537"
$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("D39").Value = "1.114.16"
$ws.Range("E39").Value = "  -2.60%  "

$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").Value = "'5.572"
$ws.Range("E42").Value = "  +2.79%  "

$ws.Range("D43").Value = "'0.7965"
$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("D44").Value = "'97.59"

$ws.Range("D45").Value = "1.763.95"

$ws.Range("E46").Value = "  -8.65%  "

$ws.Range("D47").Value = "'0.4432"
$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").Value = "'54.71"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").Value = "'0.05128"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("D51").Value = "'7.533"
$ws.Range("E51").Value = "  +2.52%  "
